$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("istoric")

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "skiuileuf"
$ws.Range("C6").Value = "Alt comentariu de test care contine RPA"
$ws.Range("D6").Value = 46035.0402861227

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "taminandreea"
$ws.Range("C7").Value = "RPA"
$ws.Range("D7").Value = 46035.0403193634

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "skiuileuf"
$ws.Range("C8").Value = "RPA"
$ws.Range("D8").Value = 46035.0403602083

$ws.Range("D5").Copy()
$ws.Range("D6:D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Save()
